# Update "Forecast Comparison" sheet with a new Week_Start_Date column
# and corrected forecast output.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before the existing column B (ASIN). This shifts the
# existing B:I columns (ASIN .. is_holiday_week) one column to the right,
# becoming C:J, with all their values/types preserved automatically.
$ws.Range("B:B").Insert()

# New column header
$ws.Range("B1").Value = "Week_Start_Date"

# Week_Start_Date values per week (stored as plain text, not dates)
$weekStarts = @(
    "2024-12-08",
    "2024-12-15",
    "2024-12-22",
    "2024-12-29",
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23"
)

for ($i = 0; $i -lt $weekStarts.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Range("B$row")
    $cell.NumberFormat = "@"
    $cell.Value = $weekStarts[$i]
}

# Simplify the Week labels in column A: "W01".."W09" -> "W1".."W9"
# ("W10".."W16" are already unpadded and stay unchanged.)
for ($w = 1; $w -le 9; $w++) {
    $row = $w + 1
    $ws.Range("A$row").Value = "W$w"
}

# Make sure the boolean "is_holiday_week" column (now J) keeps its boolean
# type after the shift.
for ($row = 2; $row -le 17; $row++) {
    $ws.Range("J$row").Value = $false
}
